# Dream Checklist.docx - "Phase 3" section -> "Phase 4 POSTS!!!!!" section
# Revisiting sidebar due to collapsing issues.

$d = $word.ActiveDocument

function New-OpenXmlPackage($bodyXml) {
    return @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$bodyXml
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
}

# Locate the anchor paragraph ("Phase 3 ...") by its distinctive text so we
# are not relying on a hard-coded paragraph index.
$anchorText = "Phase 3 (Styling Refinement of Current Features)"
$pIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "$anchorText*") {
        $pIndex = $i
        break
    }
}

# --- Paragraph 1: "Phase 3 (...): [Bold = Done]" -> "Phase 4 POSTS!!!!!" ---
$p1 = $d.Paragraphs.Item($pIndex)
$p1Xml = @'
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr>
<w:ilvl w:val="0"/>
<w:numId w:val="1"/>
</w:numPr>
<w:rPr>
<w:b/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:b/>
</w:rPr>
<w:t xml:space="preserve">Phase </w:t>
</w:r>
<w:r>
<w:rPr>
<w:b/>
</w:rPr>
<w:t>4 POSTS!!!!!</w:t>
</w:r>
</w:p>
'@
[void]$p1.Range.InsertXML((New-OpenXmlPackage $p1Xml))

# --- Paragraph 2: "Create pages for each sidebar feed" -> "Front-End" (ilvl 0 -> 1) ---
$p2 = $d.Paragraphs.Item($pIndex + 1)
$p2Xml = @'
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr>
<w:ilvl w:val="1"/>
<w:numId w:val="1"/>
</w:numPr>
<w:rPr>
<w:b/>
</w:rPr>
</w:pPr>
<w:r>
<w:t>Front-End</w:t>
</w:r>
</w:p>
'@
[void]$p2.Range.InsertXML((New-OpenXmlPackage $p2Xml))

# --- Paragraph 3: "Try to increase code reuse" -> sidebar resizing bug (ilvl 0 -> 2) ---
$p3 = $d.Paragraphs.Item($pIndex + 2)
$p3Xml = @'
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr>
<w:ilvl w:val="2"/>
<w:numId w:val="1"/>
</w:numPr>
<w:rPr>
<w:b/>
</w:rPr>
</w:pPr>
<w:r>
<w:t>Sidebar resizing and when it collapses is wrong, look at GitHub for commented out code</w:t>
</w:r>
</w:p>
'@
[void]$p3.Range.InsertXML((New-OpenXmlPackage $p3Xml))

# --- Paragraph 4: "Bug: Links routing ..." -> images/videos resize bug (ilvl 0 -> 2) ---
$p4 = $d.Paragraphs.Item($pIndex + 3)
$p4Xml = @'
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr>
<w:ilvl w:val="2"/>
<w:numId w:val="1"/>
</w:numPr>
<w:rPr>
<w:b/>
</w:rPr>
</w:pPr>
<w:r>
<w:t>Images and videos are too big at some part of the screen resize (could be related to the sidebar thing)</w:t>
</w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
'@
[void]$p4.Range.InsertXML((New-OpenXmlPackage $p4Xml))

# --- Paragraph 5: "Clean up css files" -> removed entirely ---
$p5 = $d.Paragraphs.Item($pIndex + 4)
$p5.Range.Delete()
